$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plantonistas")
$ws.Activate()
$ws.Range("A6").Value = "mamadou e relax"
$ws.Range("A6").Select()
